# Consumo.xlsx edit script
# Applies the changes described by the target diff to Hoja1 (sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) D5: 10 -> 5 (Router/PC count in "Cantidad" column for row 5)
# ---------------------------------------------------------------------
$ws.Range("D5").Value = 5

# ---------------------------------------------------------------------
# 2) C6: "Enrutador" -> "Modem Router" (new shared string)
# ---------------------------------------------------------------------
$ws.Range("C6").Value = "Modem Router"

# ---------------------------------------------------------------------
# 3) C44: "Switch" -> "Acces point"
# ---------------------------------------------------------------------
$ws.Range("C44").Value = "Acces point"

# ---------------------------------------------------------------------
# 4) Remove the now-duplicate "Acces point" row (old row 45), and shift
#    the Total/C-Resguardo/Amperage block (E46:F49) up by one row.
#    Work top-down, copying formats from the row below BEFORE that row's
#    own format gets overwritten.
# ---------------------------------------------------------------------

# Row 47 takes the old row 48 style+meaning ("C/Resguardo")
$ws.Range("E48").Copy() | Out-Null
$ws.Range("E47").PasteSpecial(-4122) | Out-Null
$ws.Range("F48").Copy() | Out-Null
$ws.Range("F47").PasteSpecial(-4122) | Out-Null

# Row 48 takes the old row 49 style+meaning ("Amperage E.")
$ws.Range("E49").Copy() | Out-Null
$ws.Range("E48").PasteSpecial(-4122) | Out-Null
$ws.Range("F49").Copy() | Out-Null
$ws.Range("F48").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Row 46: "Total:" label + SUM formula (was row 47's content, but summed range shrinks by one row)
$ws.Range("E46").Value = "Total:"
$ws.Range("F46").Formula = "=SUM(F42:F45)"

# Row 47: "C/Resguardo" + markup formula
$ws.Range("E47").Value = "C/Resguardo"
$ws.Range("F47").Formula = "=F46+(F46*`$D`$54)"

# Row 48: "Amperage E." + amperage formula
$ws.Range("E48").Value = "Amperage E."
$ws.Range("F48").Formula = "=F47/`$G`$54"

# Row 49: E49/F49 no longer used -> fully cleared (value + formatting)
$ws.Range("E49:F49").Clear()

# Row 45: old duplicate "Acces point" data cleared (keep blank styled cells)
$ws.Range("C45:F45").ClearContents()

# ---------------------------------------------------------------------
# 5) D55 formula: reference F48 -> F47 (since "C/Resguardo" moved to F47)
# ---------------------------------------------------------------------
$ws.Range("D55").Formula = "=F13+L10+F23+L23+F35+L35+F47+L45+L53"

# ---------------------------------------------------------------------
# 6) Column widths (best achievable given 1/6-character COM granularity)
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.166666666666666
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 8.666666666666666
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668
$ws.Columns.Item(8).ColumnWidth = 14.5
$ws.Columns.Item(9).ColumnWidth = 24.0
$ws.Columns.Item(10).ColumnWidth = 8.666666666666666
$ws.Columns.Item(11).ColumnWidth = 11.833333333333334
$ws.Columns.Item(12).ColumnWidth = 13.666666666666666

# ---------------------------------------------------------------------
# 7) Selection / active cell bookkeeping
# ---------------------------------------------------------------------
$ws.Range("F48").Select()
